$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.192232966423035
$ws.Range("B1").Value = 2.584609508514404
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 2.196079730987549
$ws.Range("E1").Value = 1.179997563362122
